$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "END_ROW" template marker column (J) used to control
# per-row termination behaviour for the ForEach/Continue/EndLoop block.
$ws.Range("J2").Value = "#! END_ROW true"
$ws.Range("J3").Value = "#! END_ROW true"
$ws.Range("J4").Value = "#! END_ROW"

$ws.Range("J9").Select()
